$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$ws.Range("A1").Value = "test"
